$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# auto-converting to floating point numbers.
$textCells = @("D5", "D6", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D42", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.934.83"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.061.39"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "587.44"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "131.13"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.060.63"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "33.67"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.119"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.571.11"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").Value = "62.069.36"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "3.071.79"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "448.06"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "80.63"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "2.59"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "7.45"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").Value = "6.49"
$ws.Range("E31").Value = "  -5.80%  "
$ws.Range("D32").Value = "25.99"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").Value = "0.0974"
$ws.Range("E33").Value = "  -5.71%  "
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "0.979"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").Value = "5.70"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").Value = "50.57"
$ws.Range("D38").Value = "0.0₃0693"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "0.0377"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").Value = "7.94"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").Value = "380.54"
$ws.Range("E42").Value = "  -7.46%  "
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").Value = "2.700.86"
$ws.Range("E44").Value = "  -7.02%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "125.16"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "0.241"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "34.66"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "2.02"
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "24.09"
$ws.Range("E51").Value = "  -4.62%  "
